$d = $word.ActiveDocument

$replacements = @(
    @{old="165÷5=33, 0"; new="429÷6=71, 3"},
    @{old="111÷3=37, 0"; new="890÷9=98, 8"},
    @{old="155÷9=17, 2"; new="812÷3=270, 2"},
    @{old="848÷9=94, 2"; new="960÷7=137, 1"},
    @{old="514÷2=257, 0"; new="877÷8=109, 5"},
    @{old="180÷5=36, 0"; new="301÷2=150, 1"},
    @{old="123÷7=17, 4"; new="108÷9=12, 0"},
    @{old="672÷8=84, 0"; new="782÷3=260, 2"},
    @{old="760÷8=95, 0"; new="688÷9=76, 4"},
    @{old="467÷8=58, 3"; new="581÷4=145, 1"},
    @{old="411÷7=58, 5"; new="328÷7=46, 6"},
    @{old="979÷8=122, 3"; new="674÷6=112, 2"},
    @{old="760÷6=126, 4"; new="693÷6=115, 3"},
    @{old="816÷4=204, 0"; new="449÷7=64, 1"},
    @{old="507÷6=84, 3"; new="889÷7=127, 0"},
    @{old="816÷7=116, 4"; new="581÷7=83, 0"},
    @{old="958÷9=106, 4"; new="975÷9=108, 3"},
    @{old="335÷6=55, 5"; new="509÷3=169, 2"},
    @{old="341÷5=68, 1"; new="335÷2=167, 1"},
    @{old="144÷5=28, 4"; new="972÷5=194, 2"},
    @{old="994÷9=110, 4"; new="830÷2=415, 0"},
    @{old="188÷6=31, 2"; new="830÷4=207, 2"},
    @{old="647÷2=323, 1"; new="751÷4=187, 3"},
    @{old="134÷2=67, 0"; new="494÷7=70, 4"},
    @{old="372÷7=53, 1"; new="915÷6=152, 3"}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
